$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "25.750.58"
$ws.Range("E2").Value = "  -0.14%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.631.34"
$ws.Range("E3").Value = "  -0.25%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.14%  "

# Row 5 - BNB
$ws.Range("D5").Value = "215.21"
$ws.Range("E5").Value = "  -0.01%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.81%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.14%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.64%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -1.44%  "

# Row 10 - Solana
$ws.Range("D10").Value = "19.48"
$ws.Range("E10").Value = "  -1.73%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.89%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "1.857.62"
$ws.Range("E13").Value = "  -0.19%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.631.50"
$ws.Range("E14").Value = "  -0.29%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.27%  "

# Row 16 - ShibaInu
$ws.Range("D16").Value = "0.0" + [char]0x2083 + "0762"
$ws.Range("E16").Value = "  -1.82%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "63.04"
$ws.Range("E17").Value = "  +0.03%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "25.764.58"
$ws.Range("E18").Value = "  -0.18%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.11%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  -0.14%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "192.10"
$ws.Range("E21").Value = "  -0.85%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "9.93"
$ws.Range("E22").Value = "  +0.06%  "

# Row 23 - Chainlink
$ws.Range("E23").Value = "  +1.82%  "

# Row 24 - BinanceUSD
$ws.Range("E24").Value = "  -0.07%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +2.69%  "

# Row 26 - Monero
$ws.Range("D26").Value = "143.02"
$ws.Range("E26").Value = "  +2.42%  "

# Row 28 - Cosmos
$ws.Range("D28").Value = "6.86"
$ws.Range("E28").Value = "  +0.52%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "15.48"
$ws.Range("E29").Value = "  -0.41%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.12%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -0.83%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  +0.10%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "3.22"
$ws.Range("E33").Value = "  -0.69%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  -1.78%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  -0.51%  "

# Row 36 - ARBITRUM
$ws.Range("D36").Value = "0.902"

# Row 37 - Maker
$ws.Range("D37").Value = "1.131.38"
$ws.Range("E37").Value = "  +2.20%  "

# Row 38 - MXToken
$ws.Range("E38").Value = "  -1.91%  "

# Row 39 - ImmutableX
$ws.Range("D39").Value = "0.542"

# Row 40 - VeChain
$ws.Range("E40").Value = "  -1.43%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.10%  "

# Row 42 - mCoin
$ws.Range("D42").Value = "2.53"
$ws.Range("E42").Value = "  +0.95%  "

# Row 43 - now Quant (was FraxShare)
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "100.83"
$ws.Range("E43").Value = "  +1.63%  "

# Row 44 - now FraxShare (was Quant)
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.55"
$ws.Range("E44").Value = "  -0.61%  "

# Row 45 - TrustWalletToken
$ws.Range("D45").Value = "0.796"
$ws.Range("E45").Value = "  -0.57%  "

# Row 46 - RocketPoolETH
$ws.Range("D46").Value = "1.765.77"

# Row 47 - BabyDogeCoin
$ws.Range("E47").Value = "  +3.51%  "

# Row 48 - Aave
$ws.Range("D48").Value = "55.31"
$ws.Range("E48").Value = "  -0.56%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  +0.74%  "

# Row 51 - SynthetixNetwork
$ws.Range("D51").Value = "2.34"
$ws.Range("E51").Value = "  -8.27%  "
